$d = $word.ActiveDocument

# Anchor the edit on the "Objective 1: Results" paragraph rather than a
# hard-coded index.
$find = $d.Content.Find
[void]$find.Execute("Objective 1: Results", $true, $false, $false, $false,
                     $false, $true, 1, $false, "", 0)

if ($find.Found) {
    $objParaIndex = $find.Parent.Paragraphs.Item(1).Index
} else {
    # Fallback: the heading is the first paragraph in this document.
    $objParaIndex = 1
}

# Insert three new paragraphs right after "Objective 1: Results":
#   1) a blank separator line
#   2) "Seeing how data scientist make the most compared to other data
#      related jobs"
#   3) "Seeing how Attorneys make the most compared to other non-data
#      related jobs"
# This leaves the two pre-existing blank paragraphs that lead into the
# table untouched, just pushed further down.
$d.Paragraphs.Item($objParaIndex).Range.InsertParagraphAfter()
$d.Paragraphs.Item($objParaIndex + 1).Range.InsertParagraphAfter()
$d.Paragraphs.Item($objParaIndex + 2).Range.InsertParagraphAfter()

$d.Paragraphs.Item($objParaIndex + 2).Range.Text = "Seeing how data scientist make the most compared to other data related jobs"
$d.Paragraphs.Item($objParaIndex + 3).Range.Text = "Seeing how Attorneys make the most compared to other non-data related jobs"
